$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header to reflect the new price date
$ws.Range("D1").Value = "Price 2021-05-12"

# The price column (D2:D19) stores its numbers as TEXT (shared strings),
# so force a text number format before assigning, otherwise Excel will
# silently convert these numeric-looking strings into real numbers.
$priceRange = $ws.Range("D2:D19")
$priceRange.NumberFormat = "@"

$prices = @{
    2  = "249.99"
    3  = "0.35"
    4  = "20.99"
    5  = "5.99"
    6  = "39.99"
    7  = "26.99"
    8  = "22.99"
    9  = "7.99"
    10 = "7.99"
    11 = "12.99"
    12 = "3.49"
    13 = "4.49"
    14 = "8.49"
    15 = "4.49"
    16 = "10.99"
    17 = "10.99"
    18 = "89.99"
    19 = "27.99"
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 4).Value = $prices[$row]
}

# Restore the plain/default cell style so we don't leave a leftover
# explicit number-format applied to the cells themselves.
$priceRange.Style = "Normal"
